$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F6").Value = 258
$ws.Range("F7").Value = 13062
$ws.Range("F8").Value = 58
$ws.Range("F10").Value = 262
$ws.Range("F11").Value = 3106
$ws.Range("F13").Value = 6521
$ws.Range("F16").Value = 3450
$ws.Range("F17").Value = 166
$ws.Range("F18").Value = 125
$ws.Range("F22").Value = 121
$ws.Range("F24").Value = 3628
$ws.Range("F25").Value = 94
$ws.Range("F27").Value = 2873
$ws.Range("F28").Value = 414
$ws.Range("F29").Value = 1894
$ws.Range("F31").Value = 219
$ws.Range("F32").Value = 6687
$ws.Range("F34").Value = 1066
$ws.Range("F35").Value = 1992
$ws.Range("F36").Value = 1294
$ws.Range("F37").Value = 102
$ws.Range("F38").Value = 1044
$ws.Range("F40").Value = 214
$ws.Range("F41").Value = 223
$ws.Range("F42").Value = 1150
$ws.Range("F43").Value = 1143
$ws.Range("F44").Value = 139
$ws.Range("F45").Value = 1210
$ws.Range("F46").Value = 1792
$ws.Range("F47").Value = 66

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 101

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 441
$ws.Range("F3").Value = 610
$ws.Range("F4").Value = 19

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F6").Value = 441
$ws.Range("F7").Value = 610
$ws.Range("F8").Value = 258
$ws.Range("F9").Value = 13062
$ws.Range("F10").Value = 58
$ws.Range("F13").Value = 262
$ws.Range("F14").Value = 3106
$ws.Range("F16").Value = 3450
$ws.Range("F17").Value = 166
$ws.Range("F22").Value = 121
$ws.Range("F24").Value = 3628
$ws.Range("F27").Value = 2873
$ws.Range("F28").Value = 2873
$ws.Range("F29").Value = 414
$ws.Range("F30").Value = 1894
$ws.Range("F32").Value = 219
$ws.Range("F33").Value = 6687
$ws.Range("F34").Value = 101
$ws.Range("F36").Value = 1066
$ws.Range("F37").Value = 1992
$ws.Range("F39").Value = 1294
$ws.Range("F40").Value = 102
$ws.Range("F41").Value = 1044
$ws.Range("F42").Value = 214
$ws.Range("F43").Value = 223
$ws.Range("F44").Value = 1150
$ws.Range("F45").Value = 1210
$ws.Range("F47").Value = 1792
$ws.Range("F48").Value = 66
